$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the H1 title ("Play Druidess Gold Online Slot Game for Free").
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. At the end of the document, the paragraph that used to hold the
#    italic "Prompt: ..." text becomes the meta-description text, and a
#    new bold paragraph repeating the title is inserted right before it.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Druidess Gold Online Slot Game for Free</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Druidess Gold, a mystical online slot game. Play for free and enjoy 1,296 ways to win, locked wilds, and magical graphics.</w:t></w:r></w:p>"
$lastPara.Range.InsertXML($xml)
